# "Estado de Cuenta" update:
#  - Roll the worker's mora-period window forward one month (now ends at
#    2507 instead of 2505) and grow it from 12 to 13 periods.
#  - Add two brand-new workers, each with a single 2507 row.
#  - Refresh the header totals (VALOR MORA, Cant. Trabajadores, Cant. Periodos).
#  - Push the signature/footer block down to make room for the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room: shove the footer (rows 32/33) down to rows 35/36 -------
$ws.Range("32:34").Insert()

# --- Row 30 becomes the new "last" (bottom-border) data row; capture   --
# --- that formatting from the current last row (27) before it's reused.
$ws.Range("B27:J27").Copy() | Out-Null
$ws.Range("B30").PasteSpecial(-4122) | Out-Null

# --- Row 27 becomes a regular interior row; rows 28 & 29 are brand new --
# --- interior rows too -- all three get row 16's formatting.
$ws.Range("B16:J16").Copy() | Out-Null
$ws.Range("B27").PasteSpecial(-4122) | Out-Null
$ws.Range("B16:J16").Copy() | Out-Null
$ws.Range("B28").PasteSpecial(-4122) | Out-Null
$ws.Range("B16:J16").Copy() | Out-Null
$ws.Range("B29").PasteSpecial(-4122) | Out-Null

# --- JOSE MIGUEL PACHECO MOTEZUMA: 13 periods, newest (2507) first -----
$periods = @("2507","2505","2504","2503","2502","2501","2412","2411","2410","2409","2408","2407","2406")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = "1047455782"
    $ws.Cells.Item($r, 4).Value = "JOSE MIGUEL PACHECO MOTEZUMA"
    $ws.Cells.Item($r, 5).Value = $periods[$i]
    $ws.Cells.Item($r, 6).Value = 111280
    $ws.Cells.Item($r, 7).Value = 2782000
}

# --- New worker: MARLYN PAOLA MONTERO ESCOLAR (row 29) -----------------
$ws.Cells.Item(29, 2).Value = "CC"
$ws.Cells.Item(29, 3).Value = "1143336849"
$ws.Cells.Item(29, 4).Value = "MARLYN PAOLA MONTERO ESCOLAR"
$ws.Cells.Item(29, 5).Value = "2507"
$ws.Cells.Item(29, 6).Value = 80000
$ws.Cells.Item(29, 7).Value = 2000000

# --- New worker: EDGARDO DE JESUS HERNANDEZ GARCIA (row 30) ------------
$ws.Cells.Item(30, 2).Value = "CC"
$ws.Cells.Item(30, 3).Value = "1047473869"
$ws.Cells.Item(30, 4).Value = "EDGARDO DE JESUS HERNANDEZ GARCIA"
$ws.Cells.Item(30, 5).Value = "2507"
$ws.Cells.Item(30, 6).Value = 80000
$ws.Cells.Item(30, 7).Value = 2000000

# --- Header totals -------------------------------------------------------
$ws.Range("E11").Value = 1606640   # VALOR MORA
$ws.Range("C13").Value = 3         # Cant. Trabajadores
$ws.Range("F13").Value = 13        # Cant. Periodos

# --- Column D ("Nombre Trabajador") is best-fit; widen it now that it  --
# --- holds the longer new name.
$ws.Columns.Item(4).AutoFit() | Out-Null
